$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 30.07831966666667
$ws.Cells.Item(2, 8).Value = 90.234959
$ws.Cells.Item(2, 9).Value = 0.2269842729019557
$ws.Cells.Item(2, 10).Value = 0.2269842729019557
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 31.22896466666667
$ws.Cells.Item(2, 14).Value = 93.686894
$ws.Cells.Item(2, 15).Value = 0.2877106972998646
$ws.Cells.Item(2, 16).Value = 0.2877106972998646
$ws.Cells.Item(2, 17).Value = 939.3147821030385
$ws.Cells.Item(2, 18).Value = 8453.833038927345
$ws.Cells.Item(2, 19).Value = 0.06530580343272444
$ws.Cells.Item(2, 20).Value = 0.06530580343272444

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 30.07831966666667
$ws.Cells.Item(3, 8).Value = 90.234959
$ws.Cells.Item(3, 9).Value = 0.2269842729019557
$ws.Cells.Item(3, 10).Value = 0.2269842729019557
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 40.44578266666667
$ws.Cells.Item(3, 14).Value = 121.337348
$ws.Cells.Item(3, 15).Value = 0.3726247238124506
$ws.Cells.Item(3, 16).Value = 0.3726247238124505
$ws.Cells.Item(3, 17).Value = 1216.541180216526
$ws.Cells.Item(3, 18).Value = 10948.87062194873
$ws.Cells.Item(3, 19).Value = 0.08457995199986117
$ws.Cells.Item(3, 20).Value = 0.08457995199986115

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 30.07831966666667
$ws.Cells.Item(4, 8).Value = 90.234959
$ws.Cells.Item(4, 9).Value = 0.2269842729019557
$ws.Cells.Item(4, 10).Value = 0.2269842729019557
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 25.36964133333333
$ws.Cells.Item(4, 14).Value = 76.108924
$ws.Cells.Item(4, 15).Value = 0.2337290805561598
$ws.Cells.Item(4, 16).Value = 0.2337290805561598
$ws.Cells.Item(4, 17).Value = 763.0761818526796
$ws.Cells.Item(4, 18).Value = 6867.685636674117
$ws.Cells.Item(4, 19).Value = 0.05305282540608256
$ws.Cells.Item(4, 20).Value = 0.05305282540608256

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 30.07831966666667
$ws.Cells.Item(5, 8).Value = 90.234959
$ws.Cells.Item(5, 9).Value = 0.2269842729019557
$ws.Cells.Item(5, 10).Value = 0.2269842729019557
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 11.49855033333333
$ws.Cells.Item(5, 14).Value = 34.495651
$ws.Cells.Item(5, 15).Value = 0.1059354983315251
$ws.Cells.Item(5, 16).Value = 0.1059354983315251
$ws.Cells.Item(5, 17).Value = 345.8570726292566
$ws.Cells.Item(5, 18).Value = 3112.713653663309
$ws.Cells.Item(5, 19).Value = 0.02404569206328758
$ws.Cells.Item(5, 20).Value = 0.02404569206328758

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 31.96959266666667
$ws.Cells.Item(6, 8).Value = 95.90877800000001
$ws.Cells.Item(6, 9).Value = 0.2412566535243296
$ws.Cells.Item(6, 10).Value = 0.2412566535243296
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 31.22896466666667
$ws.Cells.Item(6, 14).Value = 93.686894
$ws.Cells.Item(6, 15).Value = 0.2877106972998646
$ws.Cells.Item(6, 16).Value = 0.2877106972998646
$ws.Cells.Item(6, 17).Value = 998.3772797950593
$ws.Cells.Item(6, 18).Value = 8985.395518155532
$ws.Cells.Item(6, 19).Value = 0.06941212001371672
$ws.Cells.Item(6, 20).Value = 0.06941212001371672

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 31.96959266666667
$ws.Cells.Item(7, 8).Value = 95.90877800000001
$ws.Cells.Item(7, 9).Value = 0.2412566535243296
$ws.Cells.Item(7, 10).Value = 0.2412566535243296
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 40.44578266666667
$ws.Cells.Item(7, 14).Value = 121.337348
$ws.Cells.Item(7, 15).Value = 0.3726247238124506
$ws.Cells.Item(7, 16).Value = 0.3726247238124505
$ws.Cells.Item(7, 17).Value = 1293.035196937861
$ws.Cells.Item(7, 18).Value = 11637.31677244074
$ws.Cells.Item(7, 19).Value = 0.08989819388741942
$ws.Cells.Item(7, 20).Value = 0.0898981938874194

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 31.96959266666667
$ws.Cells.Item(8, 8).Value = 95.90877800000001
$ws.Cells.Item(8, 9).Value = 0.2412566535243296
$ws.Cells.Item(8, 10).Value = 0.2412566535243296
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 25.36964133333333
$ws.Cells.Item(8, 14).Value = 76.108924
$ws.Cells.Item(8, 15).Value = 0.2337290805561598
$ws.Cells.Item(8, 16).Value = 0.2337290805561598
$ws.Cells.Item(8, 17).Value = 811.057099526097
$ws.Cells.Item(8, 18).Value = 7299.513895734873
$ws.Cells.Item(8, 19).Value = 0.05638869580629757
$ws.Cells.Item(8, 20).Value = 0.05638869580629757

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 31.96959266666667
$ws.Cells.Item(9, 8).Value = 95.90877800000001
$ws.Cells.Item(9, 9).Value = 0.2412566535243296
$ws.Cells.Item(9, 10).Value = 0.2412566535243296
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 11.49855033333333
$ws.Cells.Item(9, 14).Value = 34.495651
$ws.Cells.Item(9, 15).Value = 0.1059354983315251
$ws.Cells.Item(9, 16).Value = 0.1059354983315251
$ws.Cells.Item(9, 17).Value = 367.603970413831
$ws.Cells.Item(9, 18).Value = 3308.435733724479
$ws.Cells.Item(9, 19).Value = 0.02555764381689596
$ws.Cells.Item(9, 20).Value = 0.02555764381689596

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 55.83720666666667
$ws.Cells.Item(10, 8).Value = 167.51162
$ws.Cells.Item(10, 9).Value = 0.4213722008598541
$ws.Cells.Item(10, 10).Value = 0.4213722008598541
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 31.22896466666667
$ws.Cells.Item(10, 14).Value = 93.686894
$ws.Cells.Item(10, 15).Value = 0.2877106972998646
$ws.Cells.Item(10, 16).Value = 0.2877106972998646
$ws.Cells.Item(10, 17).Value = 1743.738154078698
$ws.Cells.Item(10, 18).Value = 15693.64338670828
$ws.Cells.Item(10, 19).Value = 0.1212332897321672
$ws.Cells.Item(10, 20).Value = 0.1212332897321672

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 55.83720666666667
$ws.Cells.Item(11, 8).Value = 167.51162
$ws.Cells.Item(11, 9).Value = 0.4213722008598541
$ws.Cells.Item(11, 10).Value = 0.4213722008598541
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 40.44578266666667
$ws.Cells.Item(11, 14).Value = 121.337348
$ws.Cells.Item(11, 15).Value = 0.3726247238124506
$ws.Cells.Item(11, 16).Value = 0.3726247238124505
$ws.Cells.Item(11, 17).Value = 2258.379525553751
$ws.Cells.Item(11, 18).Value = 20325.41572998376
$ws.Cells.Item(11, 19).Value = 0.1570136999676476
$ws.Cells.Item(11, 20).Value = 0.1570136999676476

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 55.83720666666667
$ws.Cells.Item(12, 8).Value = 167.51162
$ws.Cells.Item(12, 9).Value = 0.4213722008598541
$ws.Cells.Item(12, 10).Value = 0.4213722008598541
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 25.36964133333333
$ws.Cells.Item(12, 14).Value = 76.108924
$ws.Cells.Item(12, 15).Value = 0.2337290805561598
$ws.Cells.Item(12, 16).Value = 0.2337290805561598
$ws.Cells.Item(12, 17).Value = 1416.569906188542
$ws.Cells.Item(12, 18).Value = 12749.12915569688
$ws.Cells.Item(12, 19).Value = 0.09848693707889918
$ws.Cells.Item(12, 20).Value = 0.09848693707889918

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 55.83720666666667
$ws.Cells.Item(13, 8).Value = 167.51162
$ws.Cells.Item(13, 9).Value = 0.4213722008598541
$ws.Cells.Item(13, 10).Value = 0.4213722008598541
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 11.49855033333333
$ws.Cells.Item(13, 14).Value = 34.495651
$ws.Cells.Item(13, 15).Value = 0.1059354983315251
$ws.Cells.Item(13, 16).Value = 0.1059354983315251
$ws.Cells.Item(13, 17).Value = 642.0469313294022
$ws.Cells.Item(13, 18).Value = 5778.42238196462
$ws.Cells.Item(13, 19).Value = 0.04463827408114016
$ws.Cells.Item(13, 20).Value = 0.04463827408114016

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 14.62767266666667
$ws.Cells.Item(14, 8).Value = 43.883018
$ws.Cells.Item(14, 9).Value = 0.1103868727138606
$ws.Cells.Item(14, 10).Value = 0.1103868727138606
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 31.22896466666667
$ws.Cells.Item(14, 14).Value = 93.686894
$ws.Cells.Item(14, 15).Value = 0.2877106972998646
$ws.Cells.Item(14, 16).Value = 0.2877106972998646
$ws.Cells.Item(14, 17).Value = 456.8070728628991
$ws.Cells.Item(14, 18).Value = 4111.263655766092
$ws.Cells.Item(14, 19).Value = 0.03175948412125625
$ws.Cells.Item(14, 20).Value = 0.03175948412125625

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 14.62767266666667
$ws.Cells.Item(15, 8).Value = 43.883018
$ws.Cells.Item(15, 9).Value = 0.1103868727138606
$ws.Cells.Item(15, 10).Value = 0.1103868727138606
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 40.44578266666667
$ws.Cells.Item(15, 14).Value = 121.337348
$ws.Cells.Item(15, 15).Value = 0.3726247238124506
$ws.Cells.Item(15, 16).Value = 0.3726247238124505
$ws.Cells.Item(15, 17).Value = 591.6276695951404
$ws.Cells.Item(15, 18).Value = 5324.649026356264
$ws.Cells.Item(15, 19).Value = 0.04113287795752246
$ws.Cells.Item(15, 20).Value = 0.04113287795752246

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 14.62767266666667
$ws.Cells.Item(16, 8).Value = 43.883018
$ws.Cells.Item(16, 9).Value = 0.1103868727138606
$ws.Cells.Item(16, 10).Value = 0.1103868727138606
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 25.36964133333333
$ws.Cells.Item(16, 14).Value = 76.108924
$ws.Cells.Item(16, 15).Value = 0.2337290805561598
$ws.Cells.Item(16, 16).Value = 0.2337290805561598
$ws.Cells.Item(16, 17).Value = 371.0988090947369
$ws.Cells.Item(16, 18).Value = 3339.889281852632
$ws.Cells.Item(16, 19).Value = 0.02580062226488049
$ws.Cells.Item(16, 20).Value = 0.02580062226488049

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 14.62767266666667
$ws.Cells.Item(17, 8).Value = 43.883018
$ws.Cells.Item(17, 9).Value = 0.1103868727138606
$ws.Cells.Item(17, 10).Value = 0.1103868727138606
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 11.49855033333333
$ws.Cells.Item(17, 14).Value = 34.495651
$ws.Cells.Item(17, 15).Value = 0.1059354983315251
$ws.Cells.Item(17, 16).Value = 0.1059354983315251
$ws.Cells.Item(17, 17).Value = 168.1970304171909
$ws.Cells.Item(17, 18).Value = 1513.773273754718
$ws.Cells.Item(17, 19).Value = 0.01169388837020146
$ws.Cells.Item(17, 20).Value = 0.01169388837020146
